# Reorders rows 70-82 of the "Artfynd" sheet.
# The row contents (columns A, B, D, E, F, G, H, Q, R) get shuffled among
# the row positions 70-82 as follows (new row -> source row):
#   70<-73  71<-81  72<-77  73<-80  74<-71  75<-76  76<-82
#   77<-72  78<-70  79<-78  80<-74  81<-75  82<-79
# All other columns (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT,
# AW, AX, AY, ...) are identical across these rows, so they are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)   # A, B, D, E, F, G, H, Q, R

$firstRow = 70
$lastRow = 82

$mapNew = @(73, 81, 77, 80, 71, 76, 82, 72, 70, 78, 74, 75, 79)

# Snapshot the original values for every relevant cell before writing
# anything, since several destination rows are also source rows.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

$i = 0
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapNew[$i]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $snapshot["$srcRow-$c"]
    }
    $i = $i + 1
}
